$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the "LUX" -> "LU" country-code typo (row 31, clean column)
$ws.Range("B31").Value = "LU"

# Append the new "Unknown" / "nan" cleaning row
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "Unknown"
$ws.Range("C33").Value = "nan"

# Match the formatting already used by the rest of column A
$ws.Range("A32").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
